# Append 2024 round-24 constructor fantasy data (rows 232-241) to Sheet1,
# following the same column layout as the existing rows:
# A constructor_name | B driver_points_total | C qualifying_bonus
# D fastest_pitstop_time | E pitstop_points | F dq_penalty
# G fastest_pitstop_bonus | H world_record_bonus | I constructor_fantasy_points
# J season | K round

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(
    @{ Row = 232; Name = "Alpine F1 Team"; B = 11; C = 6;  D = $null;      E = 0;  I = 17; K = 24 },
    @{ Row = 233; Name = "Aston Martin";   B = 3;  C = 8;  D = $null;      E = 0;  I = 11; K = 24 },
    @{ Row = 234; Name = "Ferrari";        B = 58; C = 8;  D = $null;      E = 0;  I = 66; K = 24 },
    @{ Row = 235; Name = "Haas F1 Team";   B = 18; C = 8;  D = -100.411;   E = 20; I = 46; K = 24 },
    @{ Row = 236; Name = "McLaren";        B = 37; C = 13; D = $null;      E = 0;  I = 50; K = 24 },
    @{ Row = 237; Name = "Mercedes";       B = 39; C = 6;  D = $null;      E = 0;  I = 45; K = 24 },
    @{ Row = 238; Name = "RB F1 Team";     B = -6; C = 3;  D = $null;      E = 0;  I = -3; K = 24 },
    @{ Row = 239; Name = "Red Bull";       B = 13; C = 13; D = $null;      E = 0;  I = 26; K = 24 },
    @{ Row = 240; Name = "Sauber";         B = -5; C = 6;  D = $null;      E = 0;  I = 1;  K = 24 },
    @{ Row = 241; Name = "Williams";       B = 8;  C = -1; D = $null;      E = 0;  I = 7;  K = 24 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value = $r.Name      # A constructor_name
    $ws.Cells.Item($rowNum, 2).Value = $r.B          # B driver_points_total
    $ws.Cells.Item($rowNum, 3).Value = $r.C          # C qualifying_bonus

    if ($null -ne $r.D) {
        $ws.Cells.Item($rowNum, 4).Value = $r.D      # D fastest_pitstop_time
    }

    $ws.Cells.Item($rowNum, 5).Value = $r.E          # E pitstop_points
    $ws.Cells.Item($rowNum, 6).Value = 0             # F dq_penalty
    $ws.Cells.Item($rowNum, 7).Value = 0             # G fastest_pitstop_bonus
    $ws.Cells.Item($rowNum, 8).Value = 0             # H world_record_bonus
    $ws.Cells.Item($rowNum, 9).Value = $r.I          # I constructor_fantasy_points
    $ws.Cells.Item($rowNum, 10).Value = 2024         # J season
    $ws.Cells.Item($rowNum, 11).Value = $r.K         # K round
}
